# Auto-generated Excel COM-interop script to apply market-data refresh diff
# to the FFXIV "Leve Profits" workbook (Mateus_Profits.xlsx).
#
# For every changed row, the scheduled runner re-pulled current market
# prices and rewrote the derived price/profit columns (H..N) in place.
# No formulas are involved - every cell holds a literal number - so we
# simply overwrite each affected cell with its new literal value via
# Range.Value2. Two rows lose their HQ-profit (N) cell entirely (it is
# cleared instead of being set to 0), matching the target workbook state.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value2 = 5553.6055  # was 5144.6855
$ws.Range("I15").Value2 = 5553.6055  # was 5144.6855
$ws.Range("K15").Value2 = 16660.8165  # was 15434.0565
$ws.Range("M15").Value2 = -16491.8165  # was -15265.0565

$ws.Range("H135").Value2 = 1758.8125  # was 1847.0667
$ws.Range("I135").Value2 = 1836.0714  # was 1943.8462
$ws.Range("K135").Value2 = 16524.6426  # was 17494.6158
$ws.Range("M135").Value2 = -13989.6426  # was -14959.6158

$ws.Range("H137").Value2 = 3070.5715  # was 2794.9092
$ws.Range("I137").Value2 = 4000  # was 2650
$ws.Range("K137").Value2 = 12000  # was 7950
$ws.Range("M137").Value2 = -9450  # was -5400

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 2490.4583  # was 5003.722
$ws.Range("I2").Value2 = 1900.7858  # was 4946.4165
$ws.Range("J2").Value2 = 3316  # was 5118.3335
$ws.Range("K2").Value2 = 1900.7858  # was 4946.4165
$ws.Range("L2").Value2 = 3316  # was 5118.3335
$ws.Range("M2").Value2 = -1787.7858  # was -4833.4165
$ws.Range("N2").Value2 = -3542  # was -5344.3335

$ws.Range("H32").Value2 = 4708.082  # was 5021.9297
$ws.Range("I32").Value2 = 4784.9165  # was 5109.857
$ws.Range("K32").Value2 = 4784.9165  # was 5109.857
$ws.Range("M32").Value2 = -4497.9165  # was -4822.857

$ws.Range("H61").Value2 = 9268002  # was 7944289.5
$ws.Range("I61").Value2 = 16673415  # was 12826169
$ws.Range("K61").Value2 = 16673415  # was 12826169
$ws.Range("M61").Value2 = -16673203  # was -12825957

$ws.Range("H74").Value2 = 5832.364  # was 6306.241
$ws.Range("I74").Value2 = 4567.8623  # was 4915.24
$ws.Range("K74").Value2 = 4567.8623  # was 4915.24
$ws.Range("M74").Value2 = -3693.8623  # was -4041.24

$ws.Range("H77").Value2 = 5832.364  # was 6306.241
$ws.Range("I77").Value2 = 4567.8623  # was 4915.24
$ws.Range("K77").Value2 = 22839.3115  # was 24576.2
$ws.Range("M77").Value2 = -18471.3115  # was -20208.2

$ws.Range("H110").Value2 = 2696.0264  # was 3821.077
$ws.Range("I110").Value2 = 2041.9286  # was 2932.3157
$ws.Range("J110").Value2 = 4527.5  # was 6233.4287
$ws.Range("K110").Value2 = 2041.9286  # was 2932.3157
$ws.Range("L110").Value2 = 4527.5  # was 6233.4287
$ws.Range("M110").Value2 = 3.07140000000004  # was -887.3157000000001
$ws.Range("N110").Value2 = -8617.5  # was -10323.4287

$ws.Range("H116").Value2 = 2490.4583  # was 5003.722
$ws.Range("I116").Value2 = 1900.7858  # was 4946.4165
$ws.Range("J116").Value2 = 3316  # was 5118.3335
$ws.Range("K116").Value2 = 1900.7858  # was 4946.4165
$ws.Range("L116").Value2 = 3316  # was 5118.3335
$ws.Range("M116").Value2 = 393.2141999999999  # was -2652.4165
$ws.Range("N116").Value2 = -7904  # was -9706.333500000001

$ws.Range("H122").Value2 = 3962.1538  # was 4505.273
$ws.Range("I122").Value2 = 2899.25  # was 3674.4
$ws.Range("J122").Value2 = 5662.8  # was 5197.6665
$ws.Range("K122").Value2 = 8697.75  # was 11023.2
$ws.Range("L122").Value2 = 16988.4  # was 15592.9995
$ws.Range("M122").Value2 = -6247.75  # was -8573.200000000001
$ws.Range("N122").Value2 = -21888.4  # was -20492.9995

$ws.Range("H132").Value2 = 3355.1553  # was 3416.6724
$ws.Range("I132").Value2 = 2966.453  # was 3051.8462
$ws.Range("J132").Value2 = 7475.4  # was 6578.5
$ws.Range("K132").Value2 = 8899.359  # was 9155.5386
$ws.Range("L132").Value2 = 22426.2  # was 19735.5
$ws.Range("M132").Value2 = -6369.359  # was -6625.5386
$ws.Range("N132").Value2 = -27486.2  # was -24795.5

$ws.Range("H136").Value2 = 9268002  # was 7944289.5
$ws.Range("I136").Value2 = 16673415  # was 12826169
$ws.Range("K136").Value2 = 50020245  # was 38478507
$ws.Range("M136").Value2 = -50017695  # was -38475957

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 2490.4583  # was 5003.722
$ws.Range("I3").Value2 = 1900.7858  # was 4946.4165
$ws.Range("J3").Value2 = 3316  # was 5118.3335
$ws.Range("K3").Value2 = 1900.7858  # was 4946.4165
$ws.Range("L3").Value2 = 3316  # was 5118.3335
$ws.Range("M3").Value2 = -1786.7858  # was -4832.4165
$ws.Range("N3").Value2 = -3544  # was -5346.3335

$ws.Range("H86").Value2 = 11766355  # was 11766311
$ws.Range("I86").Value2 = 1717.3572  # was 1686.0667
$ws.Range("J86").Value2 = 66668000  # was 100001000
$ws.Range("K86").Value2 = 1717.3572  # was 1686.0667
$ws.Range("L86").Value2 = 66668000  # was 100001000
$ws.Range("M86").Value2 = -594.3571999999999  # was -563.0667000000001
$ws.Range("N86").Value2 = -66670246  # was -100003246

$ws.Range("H89").Value2 = 11766355  # was 11766311
$ws.Range("I89").Value2 = 1717.3572  # was 1686.0667
$ws.Range("J89").Value2 = 66668000  # was 100001000
$ws.Range("K89").Value2 = 8586.786  # was 8430.333500000001
$ws.Range("L89").Value2 = 333340000  # was 500005000
$ws.Range("M89").Value2 = -2970.786  # was -2814.333500000001
$ws.Range("N89").Value2 = -333351232  # was -500016232

$ws.Range("H94").Value2 = 3043.6924  # was 3234.8333
$ws.Range("I94").Value2 = 2682.9473  # was 2720.889
$ws.Range("J94").Value2 = 4022.8572  # was 4776.6665
$ws.Range("K94").Value2 = 2682.9473  # was 2720.889
$ws.Range("L94").Value2 = 4022.8572  # was 4776.6665
$ws.Range("M94").Value2 = -2231.9473  # was -2269.889
$ws.Range("N94").Value2 = -4924.8572  # was -5678.6665

$ws.Range("H134").Value2 = 5297.4634  # was 5328.1665
$ws.Range("I134").Value2 = 5153.946  # was 5244
$ws.Range("J134").Value2 = 6625  # was 5833.1665
$ws.Range("K134").Value2 = 15461.838  # was 15732
$ws.Range("L134").Value2 = 19875  # was 17499.4995
$ws.Range("M134").Value2 = -12926.838  # was -13197
$ws.Range("N134").Value2 = -24945  # was -22569.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 4864.52  # was 5054.8335
$ws.Range("I31").Value2 = 4741.3335  # was 4584.4614
$ws.Range("J31").Value2 = 4978.231  # was 5610.727
$ws.Range("K31").Value2 = 4741.3335  # was 4584.4614
$ws.Range("L31").Value2 = 4978.231  # was 5610.727
$ws.Range("M31").Value2 = -4446.3335  # was -4289.4614
$ws.Range("N31").Value2 = -5568.231  # was -6200.727

$ws.Range("H34").Value2 = 4864.52  # was 5054.8335
$ws.Range("I34").Value2 = 4741.3335  # was 4584.4614
$ws.Range("J34").Value2 = 4978.231  # was 5610.727
$ws.Range("K34").Value2 = 4741.3335  # was 4584.4614
$ws.Range("L34").Value2 = 4978.231  # was 5610.727
$ws.Range("M34").Value2 = -4539.3335  # was -4382.4614
$ws.Range("N34").Value2 = -5382.231  # was -6014.727

$ws.Range("H58").Value2 = 7808.615  # was 5865.8335
$ws.Range("I58").Value2 = 3785.4285  # was 2547.5833
$ws.Range("K58").Value2 = 3785.4285  # was 2547.5833
$ws.Range("M58").Value2 = -3582.4285  # was -2344.5833

$ws.Range("H92").Value2 = 0  # was 7601
$ws.Range("J92").Value2 = 0  # was 7601
$ws.Range("L92").Value2 = 0  # was 7601
$ws.Range("N92").ClearContents()  # was -12593

$ws.Range("H134").Value2 = 8358  # was 6750.7144
$ws.Range("I134").Value2 = 4346.6665  # was 4251.25
$ws.Range("J134").Value2 = 14375  # was 10083.333
$ws.Range("K134").Value2 = 13039.9995  # was 12753.75
$ws.Range("L134").Value2 = 43125  # was 30249.999
$ws.Range("M134").Value2 = -10504.9995  # was -10218.75
$ws.Range("N134").Value2 = -48195  # was -35319.999

$ws.Range("H136").Value2 = 7808.615  # was 5865.8335
$ws.Range("I136").Value2 = 3785.4285  # was 2547.5833
$ws.Range("K136").Value2 = 11356.2855  # was 7642.749899999999
$ws.Range("M136").Value2 = -8806.2855  # was -5092.749899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value2 = 46156336  # was 40002320
$ws.Range("I131").Value2 = 166667200  # was 125000550
$ws.Range("J131").Value2 = 10003078  # was 9093870
$ws.Range("K131").Value2 = 500001600  # was 375001650
$ws.Range("L131").Value2 = 30009234  # was 27281610
$ws.Range("M131").Value2 = -499996560  # was -374996610
$ws.Range("N131").Value2 = -30019314  # was -27291690

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value2 = 0  # was 3000
$ws.Range("J10").Value2 = 0  # was 3000
$ws.Range("L10").Value2 = 0  # was 3000
$ws.Range("N10").ClearContents()  # was -3338

$ws.Range("H126").Value2 = 7196.25  # was 7232.0713
$ws.Range("I126").Value2 = 7243.875  # was 7279.857
$ws.Range("J126").Value2 = 7148.625  # was 7184.2856
$ws.Range("K126").Value2 = 21731.625  # was 21839.571
$ws.Range("L126").Value2 = 21445.875  # was 21552.8568
$ws.Range("M126").Value2 = -19261.625  # was -19369.571
$ws.Range("N126").Value2 = -26385.875  # was -26492.8568

$ws.Range("H132").Value2 = 2926.5806  # was 3670.6365
$ws.Range("I132").Value2 = 2896.087  # was 3701.8667
$ws.Range("J132").Value2 = 3014.25  # was 3603.7144
$ws.Range("K132").Value2 = 8688.261  # was 11105.6001
$ws.Range("L132").Value2 = 9042.75  # was 10811.1432
$ws.Range("M132").Value2 = -6158.261  # was -8575.6001
$ws.Range("N132").Value2 = -14102.75  # was -15871.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 18333  # was 7592.375
$ws.Range("I7").Value2 = 2499.5  # was 1534.1428
$ws.Range("K7").Value2 = 2499.5  # was 1534.1428
$ws.Range("M7").Value2 = -2387.5  # was -1422.1428

$ws.Range("H68").Value2 = 10777.417  # was 11348.091
$ws.Range("I68").Value2 = 10481  # was 10536.556
$ws.Range("J68").Value2 = 11666.667  # was 15000
$ws.Range("K68").Value2 = 10481  # was 10536.556
$ws.Range("L68").Value2 = 11666.667  # was 15000
$ws.Range("M68").Value2 = -9732  # was -9787.556
$ws.Range("N68").Value2 = -13164.667  # was -16498

$ws.Range("H71").Value2 = 10777.417  # was 11348.091
$ws.Range("I71").Value2 = 10481  # was 10536.556
$ws.Range("J71").Value2 = 11666.667  # was 15000
$ws.Range("K71").Value2 = 52405  # was 52682.78
$ws.Range("L71").Value2 = 58333.335  # was 75000
$ws.Range("M71").Value2 = -48661  # was -48938.78
$ws.Range("N71").Value2 = -65821.33499999999  # was -82488

$ws.Range("H82").Value2 = 1955.3572  # was 2186.8096
$ws.Range("I82").Value2 = 2147.8948  # was 2178.7058
$ws.Range("J82").Value2 = 1548.8889  # was 2221.25
$ws.Range("K82").Value2 = 2147.8948  # was 2178.7058
$ws.Range("L82").Value2 = 1548.8889  # was 2221.25
$ws.Range("M82").Value2 = -1786.8948  # was -1817.7058
$ws.Range("N82").Value2 = -2270.8889  # was -2943.25

$ws.Range("H85").Value2 = 1955.3572  # was 2186.8096
$ws.Range("I85").Value2 = 2147.8948  # was 2178.7058
$ws.Range("J85").Value2 = 1548.8889  # was 2221.25
$ws.Range("K85").Value2 = 2147.8948  # was 2178.7058
$ws.Range("L85").Value2 = 1548.8889  # was 2221.25
$ws.Range("M85").Value2 = -899.8948  # was -930.7058000000002
$ws.Range("N85").Value2 = -4044.8889  # was -4717.25

$ws.Range("H100").Value2 = 1221882.9  # was 1284574
$ws.Range("I100").Value2 = 1391172.2  # was 1473040.8
$ws.Range("K100").Value2 = 1391172.2  # was 1473040.8
$ws.Range("M100").Value2 = -1390631.2  # was -1472499.8

$ws.Range("H126").Value2 = 18333  # was 7592.375
$ws.Range("I126").Value2 = 2499.5  # was 1534.1428
$ws.Range("K126").Value2 = 7498.5  # was 4602.428400000001
$ws.Range("M126").Value2 = -5028.5  # was -2132.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 830.2308  # was 2168.5
$ws.Range("I107").Value2 = 571.875  # was 966.6667
$ws.Range("J107").Value2 = 1243.6  # was 2889.6
$ws.Range("K107").Value2 = 1715.625  # was 2900.0001
$ws.Range("L107").Value2 = 3730.8  # was 8668.799999999999
$ws.Range("M107").Value2 = 204.375  # was -980.0001000000002
$ws.Range("N107").Value2 = -7570.799999999999  # was -12508.8

$ws.Range("H126").Value2 = 2852.2563  # was 3026.7222
$ws.Range("I126").Value2 = 3226.9678  # was 3398.2415
$ws.Range("J126").Value2 = 1400.25  # was 1487.5714
$ws.Range("K126").Value2 = 9680.903399999999  # was 10194.7245
$ws.Range("L126").Value2 = 4200.75  # was 4462.7142
$ws.Range("M126").Value2 = -7210.903399999999  # was -7724.7245
$ws.Range("N126").Value2 = -9140.75  # was -9402.7142

$ws.Range("H132").Value2 = 6396.309  # was 6543.34
$ws.Range("I132").Value2 = 5809.4863  # was 5998.6
$ws.Range("K132").Value2 = 17428.4589  # was 17995.8
$ws.Range("M132").Value2 = -14898.4589  # was -15465.8

$ws.Range("H136").Value2 = 5817.2573  # was 6024.394
$ws.Range("I136").Value2 = 5315.7812  # was 5510.2
$ws.Range("K136").Value2 = 15947.3436  # was 16530.6
$ws.Range("M136").Value2 = -13397.3436  # was -13980.6
